# Applies the English -> French translation edits described by the diff.
# Uses paragraph-scoped Find/Replace so duplicated English phrases
# elsewhere in the document are not affected by a single replacement.

$d = $word.ActiveDocument
$wdReplaceAll = 2
$wdFindContinue = 1

function Replace-InRange($range, [string]$find, [string]$replace) {
    $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replace, $wdReplaceAll) | Out-Null
}

$paras = $d.Paragraphs

# Para 1: language switcher line (hyperlink "English" + trailing language list)
Replace-InRange $paras.Item(1).Range "English" "Anglais"
Replace-InRange $paras.Item(1).Range " / Portuguese / French / Thai / Vietnamese / Spanish" " / portugais / français / thaïlandais / vietnamien / espagnol"

# Para 3: "English" heading for this translation section
Replace-InRange $paras.Item(3).Range "English" "Anglais"

# Para 5: table cell heading "Brief:"
Replace-InRange $paras.Item(5).Range "Brief" "Résumé"

# Para 6: table cell body text (brief description)
Replace-InRange $paras.Item(6).Range "An email sent to partners in the target country who RSVPed yes but haven’t sent their documents to us. It will be sent via customer.io" "Un email envoyé aux partenaires du pays cible qui ont répondu oui mais qui n'ont pas encore envoyé leurs documents. Il sera envoyé via customer.io"

# Para 8: table cell heading "Target audience:"
Replace-InRange $paras.Item(8).Range "Target audience" "Public cible"

# Para 9: table cell body text (target audience)
Replace-InRange $paras.Item(9).Range "Invited partners who haven’t submitted their documents" "Partenaires qui n'ont pas soumis leurs documents"

# Para 12: "Subject line: [EVENT NAME] — have you submitted your docs?  " (first email)
Replace-InRange $paras.Item(12).Range "Subject line" "Objet"
Replace-InRange $paras.Item(12).Range ": " " : "
Replace-InRange $paras.Item(12).Range " — have you submitted your docs?  " " — avez-vous déjà envoyé vos documents ?  "

# Para 14: "Don’t forget to send your documents" heading (first email)
Replace-InRange $paras.Item(14).Range "Don’t forget to send your documents" "N'oubliez pas d'envoyer vos documents !"

# Para 16: "Hi [PARTNER NAME], " greeting (first email)
Replace-InRange $paras.Item(16).Range "Hi " "Salut "
Replace-InRange $paras.Item(16).Range "[PARTNER NAME]" "[NOM DU PARTENAIRE]"

# Para 18: "We’re excited to see you at the upcoming [EVENT NAME]. " (first email)
Replace-InRange $paras.Item(18).Range "We’re excited to see you at the upcoming " "Nous serons ravis de vous rencontrer au prochain "

# Para 19: "To confirm your registration..." + trailing colon (first email)
Replace-InRange $paras.Item(19).Range "To confirm your registration, we need the following documents from you by " "Pour confirmer votre inscription, vous devez nous envoyer les documents suivants au plus tard le "
Replace-InRange $paras.Item(19).Range ":" " :"

# Para 21: "Please send a copy of these documents..." (first email)
Replace-InRange $paras.Item(21).Range "Please send a copy of these documents to your country manager, " "Veuillez envoyer une copie de ces documents à "
Replace-InRange $paras.Item(21).Range ", at " ", votre responsable local, à l'adresse "
Replace-InRange $paras.Item(21).Range " or " " ou au numéro "
Replace-InRange $paras.Item(21).Range " (WhatsApp), so that we can make the necessary arrangements for you, including accommodation and transportation." " (WhatsApp). Cela nous permettra de prendre les dispositions nécessaires à votre égard, notamment en matière d'hébergement et de transport."

# Para 22: "If you have any questions, please contact your country manager." (first email, no comma)
Replace-InRange $paras.Item(22).Range "If you have any questions, please contact your country manager." "Si vous avez des questions, veuillez contacter votre responsable local."

# Para 23: "We look forward to seeing you there!" (first email)
Replace-InRange $paras.Item(23).Range "We look forward to seeing you there!" "Nous sommes impatients de vous y rencontrer !"

# Para 29: "Subject line: [EVENT NAME] — have you submitted your docs?  " (second email)
Replace-InRange $paras.Item(29).Range "Subject line" "Objet"
Replace-InRange $paras.Item(29).Range ": " " : "
Replace-InRange $paras.Item(29).Range " — have you submitted your docs?  " " — avez-vous déjà envoyé vos documents ?  "

# Para 31: "Don’t forget to send your documents" heading (second email)
Replace-InRange $paras.Item(31).Range "Don’t forget to send your documents" "N'oubliez pas d'envoyer vos documents !"

# Para 33: "Dear [PARTNER NAME], " greeting (second email) - only "Dear " changes
Replace-InRange $paras.Item(33).Range "Dear " "Cher "

# Para 35: "We’re excited to see you at the upcoming [EVENT NAME]. " (second email)
Replace-InRange $paras.Item(35).Range "We’re excited to see you at the upcoming " "Nous avons hâte de vous rencontrer au prochain "

# Para 36: "To ensure you have the best experience..." + trailing colon (second email)
Replace-InRange $paras.Item(36).Range "To ensure you have the best experience at this event, we need the following documents from you by " "Nous souhaitons vous offrir la meilleure expérience possible lors de cet événement. Pour cela, nous aurons besoin des documents suivants au plus tard le "
Replace-InRange $paras.Item(36).Range ":" " :"

# Para 38: "Please reply to this email with a copy of these documents..." (second email)
Replace-InRange $paras.Item(38).Range "Please reply to this email with a copy of these documents so that we have make the necessary arrangements for you, including accommodation and transportation." "Veuillez répondre à cet e-mail en joignant une copie de ces documents afin que nous puissions prendre les dispositions nécessaires à votre égard, notamment en matière d'hébergement et de transport."

# Para 39: "If you have any questions, please contact us via live chat or WhatsApp. " (second email)
Replace-InRange $paras.Item(39).Range "If you have any questions, please contact us via " "Si vous avez des questions, veuillez nous contacter par "
Replace-InRange $paras.Item(39).Range "live chat" "chat en direct"
Replace-InRange $paras.Item(39).Range " or " " ou par "

# Para 40: "If you have any questions, please contact your country manager, ..." (second email, with comma)
Replace-InRange $paras.Item(40).Range "If you have any questions, please contact your country manager, " "Si vous avez des questions, veuillez contacter votre responsable local, "
Replace-InRange $paras.Item(40).Range ", at " ", à l'adresse "
Replace-InRange $paras.Item(40).Range " or " " ou au numéro "

# Para 41: "We look forward to seeing you there!" (second email)
Replace-InRange $paras.Item(41).Range "We look forward to seeing you there!" "Nous sommes impatients de vous y rencontrer !"

# Comment text: "choose either one" -> "choisissez l'un ou l'autre"
$comment = $d.Comments.Item(1)
Replace-InRange $comment.Range "choose either one" "choisissez l'un ou l'autre"
